$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new diary rows
$ws.Range("A5").Value = "19/2-2018"
$ws.Range("B5").Value = "Image Serching and spritesheat building"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 30

$ws.Range("A6").Value = "22/2-2018"
$ws.Range("B6").Value = "MovingObject made it drawable"
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 0

$ws.Range("A7").Value = "23/2-2018"
$ws.Range("B7").Value = "Added external windowSize"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0

# Widen column B to fit the new, longer task descriptions (target stored
# width is 37.109375 chars; the host quantizes ColumnWidth to whole
# pixels at a 6px max-digit-width, so 36.3 is the input that lands on the
# closest representable stored width).
$ws.Columns.Item(2).ColumnWidth = 36.3

# Update the selection to mirror the saved view state
$ws.Range("E9").Select()
